# Update the "丽水-漫展信息" workbook:
#  - Sheets "展览" (index 1) and "全部类型" (index 4) each gain a new
#    4th data row (rows shift from A1:J4 to A1:J5), and the event data
#    for the existing rows is reshuffled / updated with new details.
#
# Final row order (row 1 is the header, already correct):
#   Row 2: id=1  2024.02.07  丽水·新年动漫狂欢盛典   @ 飞达国际大酒店 ...
#   Row 3: id=2  2024.02.07  龙泉·崩X铁X原ONLY       @ 龙泉金沙温泉酒店 ...
#   Row 4: id=3  2024.02.14  丽水·YA●怀旧only        @ 丽水体育中心 ...
#   Row 5: id=4  2024.02.18  丽水·LPJ 现实X次元动漫展 @ 飞达国际大酒店 ...
#
# IMPORTANT — quirks of this interpreter discovered while testing:
#  1) Command-call arguments that are parenthesized expressions (e.g.
#     `Foo $a ($i + 1)`) are mis-parsed (silently dropped / argument
#     order scrambled). Every value passed to a function is therefore
#     first assigned to a plain named variable, then passed by name.
#  2) Only *function parameters* get an isolated scope. Plain
#     body-assigned variables (`$x = ...`) are NOT function-local: they
#     live in a scope shared across the whole call stack, so a callee
#     reusing a variable name clobbers the same-named variable in every
#     caller still on the stack (incl. enclosing `for`/`foreach` loop
#     counters!). To stay safe, every plain variable in this script uses
#     a name that is unique across the *entire* file.

function Set-TextCell($stcWs, $stcRow, $stcCol, $stcText) {
    # Plain string assignment. Safe for columns whose content never
    # looks like a pure number/date (names, addresses, URLs, the mixed
    # "start-end" time range in column E) -- Excel's normal cell-input
    # parsing leaves these as text, so no explicit number format is
    # needed and the cell keeps the sheet's default style.
    $stcCell = $stcWs.Cells.Item($stcRow, $stcCol)
    $stcCell.Value = $stcText
}

function Set-ForcedTextCell($ftcWs, $ftcRow, $ftcCol, $ftcText) {
    # Same as Set-TextCell, but first forces a Text number format.
    # Needed for columns whose content is a *pure* number or date
    # string (e.g. "2024.02.07" or "45") which Excel's input parsing
    # would otherwise silently coerce into a real date/number value.
    $ftcCell = $ftcWs.Cells.Item($ftcRow, $ftcCol)
    $ftcCell.NumberFormat = "@"
    $ftcCell.Value = $ftcText
}

function Set-NumberCell($sncWs, $sncRow, $sncCol, $sncNum) {
    $sncCell = $sncWs.Cells.Item($sncRow, $sncCol)
    # Columns that receive numbers are already in the default "General"
    # number format in the source sheet, so there is no need to (and we
    # deliberately avoid) touching NumberFormat here -- forcing it would
    # register a redundant custom numFmt entry in styles.xml.
    $sncCell.Value = $sncNum
}

function Set-BoolCell($sbcWs, $sbcRow, $sbcCol, $sbcVal) {
    $sbcCell = $sbcWs.Cells.Item($sbcRow, $sbcCol)
    $sbcCell.Value = $sbcVal
}

function Format-IdCell($ficWs, $ficRow) {
    # Mirrors the existing style applied to column A data cells: bold,
    # centered/top aligned, thin box border all around (style index "1"
    # in the original sheet).
    $ficCol = 1
    $ficCell = $ficWs.Cells.Item($ficRow, $ficCol)
    $ficCell.Font.Bold = $true
    $ficHAlign = -4108   # xlCenter
    $ficVAlign = -4160   # xlTop
    $ficCell.HorizontalAlignment = $ficHAlign
    $ficCell.VerticalAlignment = $ficVAlign
    $ficLineStyle = 1    # xlContinuous
    $ficWeight = 2       # xlThin
    $ficLeft = 7
    $ficTop = 8
    $ficBottom = 9
    $ficRight = 10
    $ficCell.Borders.Item($ficLeft).LineStyle = $ficLineStyle
    $ficCell.Borders.Item($ficLeft).Weight = $ficWeight
    $ficCell.Borders.Item($ficTop).LineStyle = $ficLineStyle
    $ficCell.Borders.Item($ficTop).Weight = $ficWeight
    $ficCell.Borders.Item($ficBottom).LineStyle = $ficLineStyle
    $ficCell.Borders.Item($ficBottom).Weight = $ficWeight
    $ficCell.Borders.Item($ficRight).LineStyle = $ficLineStyle
    $ficCell.Borders.Item($ficRight).Weight = $ficWeight
}

function Set-EventRow($serWs, $serRow, $serData) {
    $serColA = 1
    $serColB = 2
    $serColC = 3
    $serColD = 4
    $serColE = 5
    $serColF = 6
    $serColG = 7
    $serColH = 8
    $serColI = 9
    $serColJ = 10

    $serValA = $serData.A
    $serValB = $serData.B
    $serValC = $serData.C
    $serValD = $serData.D
    $serValE = $serData.E
    $serValF = $serData.F
    $serValG = $serData.G
    $serValH = $serData.H
    $serValI = $serData.I
    $serValJ = $serData.J

    Set-NumberCell $serWs $serRow $serColA $serValA
    Format-IdCell $serWs $serRow
    Set-ForcedTextCell $serWs $serRow $serColB $serValB
    Set-TextCell $serWs $serRow $serColC $serValC
    Set-TextCell $serWs $serRow $serColD $serValD
    Set-TextCell $serWs $serRow $serColE $serValE
    Set-NumberCell $serWs $serRow $serColF $serValF
    Set-ForcedTextCell $serWs $serRow $serColG $serValG
    Set-BoolCell $serWs $serRow $serColH $serValH
    Set-TextCell $serWs $serRow $serColI $serValI
    Set-TextCell $serWs $serRow $serColJ $serValJ
}

$mainRows = @(
    @{
        A = 1
        B = "2024.02.07"
        C = "丽水·新年动漫狂欢盛典"
        D = "中东路848号(解放街交汇) 飞达国际大酒店"
        E = "2024.02.07 09:00-02.07 17:00"
        F = 267
        G = "45"
        H = $false
        I = "https://show.bilibili.com/platform/detail.html?id=78294&msource=Msearch_colligation"
        J = "//i1.hdslb.com/bfs/openplatform/202311/lP5IkqWn1699431829470.jpeg"
    },
    @{
        A = 2
        B = "2024.02.07"
        C = "龙泉·崩X铁X原ONLY"
        D = "金沙路26-1号 龙泉金沙温泉酒店"
        E = "2024.02.07 10:30-02.07 16:30"
        F = 17
        G = "50"
        H = $false
        I = "https://show.bilibili.com/platform/detail.html?id=80714&msource=Msearch_colligation"
        J = "//i2.hdslb.com/bfs/openplatform/202401/rTvQio211704877379770.jpeg"
    },
    @{
        A = 3
        B = "2024.02.14"
        C = "丽水·YA●怀旧only"
        D = "人民街567号 丽水体育中心"
        E = "2024.02.14 09:00-02.14 17:00"
        F = 11
        G = "35"
        H = $false
        I = "https://show.bilibili.com/platform/detail.html?id=81032&msource=Msearch_colligation"
        J = "//i0.hdslb.com/bfs/openplatform/202401/LbqTNkvq1705561884633.png"
    },
    @{
        A = 4
        B = "2024.02.18"
        C = "丽水·LPJ 现实X次元动漫展"
        D = "中东路848号(解放街交汇) 飞达国际大酒店"
        E = "2024.02.18 10:00-02.18 17:00"
        F = 258
        G = "45"
        H = $false
        I = "https://show.bilibili.com/platform/detail.html?id=79437&msource=Msearch_colligation"
        J = "//i1.hdslb.com/bfs/openplatform/202312/ee5hLUN61702276208812.jpeg"
    }
)

$mainWb = $excel.ActiveWorkbook
$mainSheetIndexes = @(1, 4)

foreach ($mainSheetIndex in $mainSheetIndexes) {
    $mainWs = $mainWb.Worksheets.Item($mainSheetIndex)

    # Grow the sheet from 4 data rows to 5 data rows: insert a brand new
    # blank row 5 below the current last data row (row 4). Every data
    # row (2-5) is then (re)written below with its final target content,
    # so it does not matter that row 5 starts out empty.
    $mainNewRowNum = 5
    $mainWs.Rows.Item($mainNewRowNum).Insert()

    $mainRowCount = $mainRows.Count
    for ($mainIdx = 0; $mainIdx -lt $mainRowCount; $mainIdx++) {
        $mainTargetRow = $mainIdx + 2
        $mainRowData = $mainRows[$mainIdx]
        Set-EventRow $mainWs $mainTargetRow $mainRowData
    }
}

Write-Output "edit complete"
